$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the contents of the Range matched by Find.Execute with a
# sequence of raw WordprocessingML runs (passed as an OOXML fragment string
# containing <w:r>/<w:proofErr/> elements). Mirrors what Word itself does
# when it retypes / re-spell-checks a run and breaks it into several runs
# around words flagged by the proofing engine.
# ---------------------------------------------------------------------------
function Replace-RangeWithRuns($searchText, $runsXml) {

    $rng = $d.Content
    $find = $rng.Find
    $find.ClearFormatting()
    $find.Text = $searchText
    $find.MatchCase = $true
    $found = $find.Execute()
    if (-not $found) {
        throw "Could not find text: $searchText"
    }

    # Re-wrap the matched span in a brand-new Range object: InsertXML behaves
    # reliably on a plain Range(start, end), but leaves stray leftover text
    # when invoked directly on the Range that still belongs to a live Find.
    $target = $d.Range($rng.Start, $rng.End)

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runsXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($pkg) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. "Prapti Jigneshbhai Patel (20BCA019)" -> split around "Jigneshbhai"
# ---------------------------------------------------------------------------
$runs1 = '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">Prapti </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Jigneshbhai</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> Patel (20BCA019)</w:t></w:r>'
Replace-RangeWithRuns "Prapti Jigneshbhai Patel (20BCA019)" $runs1

# ---------------------------------------------------------------------------
# 2. "Atmik Maheshbhai Virani (20BCA021)" -> split around "Atmik"
# ---------------------------------------------------------------------------
$runs2 = '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Atmik</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> Maheshbhai Virani (20BCA021)</w:t></w:r>'
Replace-RangeWithRuns "Atmik Maheshbhai Virani (20BCA021)" $runs2

# ---------------------------------------------------------------------------
# 3. ". Chandaben Mohanbhai Patel Institute of Computer Applications " ->
#    split around "Chandaben" and "Mohanbhai"
# ---------------------------------------------------------------------------
$runs3 = '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Chandaben</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Mohanbhai</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> Patel Institute of Computer Applications </w:t></w:r>'
Replace-RangeWithRuns ". Chandaben Mohanbhai Patel Institute of Computer Applications " $runs3

# ---------------------------------------------------------------------------
# 4. "April" -> "May" (only this run's text changes; the sibling runs
#    "/","2022-2023"," of Submission ", <w:cr/> must stay exactly as they
#    were, so the whole paragraph's content is rebuilt with the neighbours
#    copied verbatim and only the first run's text swapped).
# ---------------------------------------------------------------------------
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$find.Text = "April/2022-2023 of Submission "
$find.MatchCase = $true
$found = $find.Execute()
if (-not $found) {
    throw "Could not find the April/2022-2023 paragraph text"
}
$paraTarget = $d.Range($rng.Start, $rng.End + 1)

$runs4 = '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>May</w:t></w:r>' +
         '<w:r w:rsidR="000005C3" w:rsidRPr="0025385F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>/</w:t></w:r>' +
         '<w:r w:rsidRPr="0025385F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>2022-2023</w:t></w:r>' +
         '<w:r w:rsidR="00A57C5E" w:rsidRPr="0025385F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> of Submission </w:t></w:r>' +
         '<w:r w:rsidR="00A57C5E" w:rsidRPr="0025385F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:cr/></w:r>'

$pkg4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $runs4 + '</w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

$paraTarget.InsertXML($pkg4) | Out-Null

# ---------------------------------------------------------------------------
# 5. The built-in "Normal Table" style loses its Quick-Style-gallery flag
#    (<w:qFormat/> removed from word/styles.xml).
# ---------------------------------------------------------------------------
$tableNormalStyle = $d.Styles("Normal Table")
$tableNormalStyle.QuickStyle = $false
